$d = $word.ActiveDocument

function Find-ParaIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Locate the "Чи є ... відповідь" paragraph.
# ---------------------------------------------------------------------------
$qIndex = Find-ParaIndex $d "Обгрунтуйте"

# ---------------------------------------------------------------------------
# Step 1: merge the empty paragraph that precedes it into it (deleting the
# empty paragraph's mark joins the two).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($qIndex - 1).Range.Delete()
$qIndex = $qIndex - 1

# ---------------------------------------------------------------------------
# Step 2: fix the typo Обгрунтуйте -> Обґрунтуйте (г -> ґ)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Обгрунтуйте", $true, $false, $false, $false, $false, $true, 1, $false, "Обґрунтуйте", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: italicise the whole (merged) question paragraph, paragraph mark
# included, so both runs and the paragraph mark pick up i/iCs.
# ---------------------------------------------------------------------------
$qPara = $d.Paragraphs.Item($qIndex)
$qPara.Range.Font.Italic = 1
$qPara.Range.Font.ItalicBi = 1

# ---------------------------------------------------------------------------
# Step 4: append the four new explanatory paragraphs right after it.
# Inserting a paragraph break at a range collapsed to qPara's End adds the
# new (empty) paragraph *after* qPara, instead of splitting qPara itself.
# ---------------------------------------------------------------------------
$texts = @(
    "Щоб відповісти на це питання, потрібно розглянути ключові ідеї кожного філософа і визначити, чи є між ними протиріччя.",
    "Рене Декарт вірив у розділення розуму і тіла, стверджуючи, що розум - це нефізична субстанція, яка здатна розуміти абстрактні поняття. Френсіс Бекон, з іншого боку, наголошував на важливості емпіричного спостереження та експерименту у здобутті знань. Хоча їхні підходи до розуміння когнітивних процесів у людини відрізняються, вони не обов'язково є суперечливими.",
    "Барух Спіноза вважав, що всесвіт - це єдина, нескінченна субстанція, яка керується детермінованими законами, тоді як Готфрід Лейбніц стверджував, що всесвіт складається з окремих субстанцій (монад), які пов'язані між собою у наперед встановленій гармонії. Хоча їхні погляди на природу Всесвіту відрізняються, вони не обов'язково є суперечливими.",
    "Тому можна стверджувати, що теорії Декарта і Бекона та Спінози і Лейбніца не обов'язково є суперечливими. Кожен філософ мав власний унікальний погляд на світ і пізнавальні процеси людини, і їхні ідеї можна розглядати як взаємодоповнюючі, а не суперечливі."
)

$curIndex = $qIndex
foreach ($txt in $texts) {
    $endPos = $d.Paragraphs.Item($curIndex).Range.End
    $d.Range($endPos, $endPos).InsertParagraphAfter()
    $curIndex = $curIndex + 1
    $d.Paragraphs.Item($curIndex).Range.Text = $txt
}

Write-Host "done"
